# This script updates the "cryptos" worksheet with refreshed market data
# (price and 1h volume/percentage change columns), matching the upstream
# GitHub Actions data-refresh commit. Row 23/24 and 36/37 also swap their
# coin identity (Coin name + Link) because the source ranking reordered
# those two pairs of coins.
#
# Numeric-looking price strings (e.g. "9.91", "0.999") must be written as
# TEXT (not numbers) so that formatting such as trailing zeros and
# thousand-separator dots ("67.246.24") is preserved exactly as in the
# source data. Setting .Value directly on a numeric-looking string causes
# Excel to coerce it into a real number, losing the original formatting.
# To avoid that we briefly force the cell's number format to Text ("@")
# before assigning the value, then call ClearFormats() to drop the
# temporary formatting so the cell is left with no explicit style (as in
# the original workbook) while keeping the value stored as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

function Set-PlainValue($cellRef, $value) {
    $ws.Range($cellRef).Value = $value
}


Set-PlainValue "D2" '67.246.24'
Set-PlainValue "E2" '  -1.02%  '
Set-PlainValue "D3" '3.583.76'
Set-PlainValue "E3" '  -2.01%  '
Set-TextValue "D4" '0.999'
Set-PlainValue "E4" '  -0.13%  '
Set-TextValue "D5" '575.88'
Set-PlainValue "E5" '  -3.46%  '
Set-TextValue "D6" '192.38'
Set-PlainValue "E6" '  +1.23%  '
Set-PlainValue "D7" '3.580.53'
Set-PlainValue "E7" '  -2.07%  '
Set-TextValue "D8" '0.618'
Set-PlainValue "E8" '  -0.16%  '
Set-PlainValue "E9" '  -0.04%  '
Set-TextValue "D10" '0.679'
Set-PlainValue "E10" '  -2.84%  '
Set-TextValue "D11" '0.151'
Set-PlainValue "E11" '  -1.19%  '
Set-PlainValue "E12" '  -3.71%  '
Set-TextValue "D13" '0.0000274'
Set-PlainValue "E13" '  +0.19%  '
Set-TextValue "D14" '9.91'
Set-PlainValue "E14" '  -2.47%  '
Set-PlainValue "D15" '4.148.59'
Set-PlainValue "E15" '  -2.17%  '
Set-PlainValue "D16" '3.579.98'
Set-PlainValue "E16" '  -2.21%  '
Set-PlainValue "E17" '  -1.12%  '
Set-TextValue "D18" '12.30'
Set-PlainValue "E18" '  -0.95%  '
Set-PlainValue "D19" '67.109.08'
Set-PlainValue "E19" '  -0.93%  '
Set-TextValue "D20" '18.34'
Set-PlainValue "E20" '  -2.71%  '
Set-TextValue "D21" '1.07'
Set-PlainValue "E21" '  -3.84%  '
Set-TextValue "D22" '404.30'
Set-PlainValue "E22" '  +1.07%  '
Set-PlainValue "B23" 'PancakeSwap'
Set-PlainValue "C23" 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D23" '4.21'
Set-PlainValue "E23" '  -4.48%  '
Set-PlainValue "B24" 'RenderToken'
Set-PlainValue "C24" 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D24" '12.58'
Set-PlainValue "E24" '  +12.23%  '
Set-TextValue "D25" '85.86'
Set-PlainValue "E25" '  -1.96%  '
Set-TextValue "D26" '2.93'
Set-PlainValue "E26" '  -0.69%  '
Set-TextValue "D27" '12.54'
Set-PlainValue "E27" '  +0.95%  '
Set-TextValue "D28" '6.11'
Set-PlainValue "E28" '  +0.96%  '
Set-TextValue "D29" '3.78'
Set-PlainValue "E29" '  -0.36%  '
Set-TextValue "D30" '7.98'
Set-PlainValue "E30" '  +8.74%  '
Set-TextValue "D31" '9.10'
Set-PlainValue "E31" '  -1.60%  '
Set-TextValue "D32" '31.31'
Set-PlainValue "E32" '  -1.50%  '
Set-TextValue "D33" '662.01'
Set-PlainValue "E33" '  +8.66%  '
Set-TextValue "D34" '12.18'
Set-PlainValue "E34" '  -0.97%  '
Set-PlainValue "E35" '  -1.13%  '
Set-PlainValue "B36" 'InjectiveProtocol'
Set-PlainValue "C36" 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D36" '42.93'
Set-PlainValue "E36" '  -4.25%  '
Set-PlainValue "B37" 'OKB'
Set-PlainValue "C37" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D37" '63.72'
Set-PlainValue "E37" '  -3.94%  '
Set-TextValue "D38" '0.415'
Set-PlainValue "E38" '  +5.54%  '
Set-PlainValue "E39" '  +0.16%  '
Set-PlainValue "D40" '0.0₃0781'
Set-PlainValue "E40" '  +0.90%  '
Set-TextValue "D41" '3.13'
Set-PlainValue "E41" '  +8.31%  '
Set-TextValue "D42" '2.82'
Set-PlainValue "E42" '  +10.83%  '
Set-PlainValue "D43" '3.144.53'
Set-PlainValue "E43" '  +12.52%  '
Set-TextValue "D44" '0.133'
Set-PlainValue "E44" '  -0.73%  '
Set-TextValue "D45" '0.998'
Set-PlainValue "E45" '  -0.29%  '
Set-TextValue "D46" '0.0417'
Set-PlainValue "E46" '  -1.88%  '
Set-PlainValue "E47" '  -0.65%  '
Set-PlainValue "E48" '  -3.19%  '
Set-TextValue "D49" '143.17'
Set-PlainValue "E49" '  +0.23%  '
Set-TextValue "D50" '8.60'
Set-PlainValue "E50" '  -2.51%  '
Set-TextValue "D51" '2.52'
Set-PlainValue "E51" '  -3.18%  '
